# Update the "Generate Report for Handback" timestamps.
# The handback report re-runs, so the handoff/handback datetime stamps
# for the 93c7890f-... row advance for each locale sheet, and the
# "Latest HO Xliff Generate Date" on the Overview sheet (which shares
# its value with the de-de sheet's Handoff Datetime cell) advances too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G3 - Latest HO Xliff Generate Date for 93c7890f-...md
$wsOverview.Range("G3").Value = "2016-08-16 06:40:43"

# zh-cn!H3 - Correspond Handoff Datetime for 93c7890f-... row
$wsZhCn.Range("H3").Value = "2016-08-16 06:40:38"

# zh-cn!K3 - Correspond Handback DateTime for 93c7890f-... row
$wsZhCn.Range("K3").Value = "2016-08-16 06:40:56"

# de-de!H3 - Correspond Handoff Datetime for 93c7890f-... row
# (shares the shared-string value with Overview!G3)
$wsDeDe.Range("H3").Value = "2016-08-16 06:40:43"

# de-de!K3 - Correspond Handback DateTime for 93c7890f-... row
$wsDeDe.Range("K3").Value = "2016-08-16 06:41:07"
